$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- "m" column header is actually the 3rd select bit -> rename to "s3" ---
$ws.Range("A5").Value = "s3"

# --- fixed CE (clock enable) throughout simulation: the E column in the ---
# --- truth table had CE asserted/cleared on the wrong rows; correct it ---
# CE=1 now on rows 6,8,10,12,14,16,18,20 ; CE cleared (blank) elsewhere
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = ""
$ws.Range("E8").Value = 1
$ws.Range("E9").Value = ""
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = ""
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = ""
$ws.Range("E14").Value = 1
$ws.Range("E15").Value = ""
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = ""
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = ""

# --- update the saved window/view state (scroll position + selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("O20").Select()
